$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style (G1) to new header cell H1, then set its value
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("H1").Value = "Save"

# Fill in the Save column (1 = saved, 0 = not saved) for each data row
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("H6").Value = 1
$ws.Range("H7").Value = 1
$ws.Range("H8").Value = 0
$ws.Range("H9").Value = 0
$ws.Range("H10").Value = 0
$ws.Range("H11").Value = 0
$ws.Range("H12").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("H14").Value = 0
$ws.Range("H15").Value = 0
$ws.Range("H16").Value = 0
$ws.Range("H17").Value = 0
$ws.Range("H18").Value = 0
$ws.Range("H19").Value = 0
$ws.Range("H20").Value = 0
$ws.Range("H21").Value = 1
$ws.Range("H22").Value = 0
$ws.Range("H23").Value = 1
$ws.Range("H24").Value = 0
$ws.Range("H25").Value = 1
$ws.Range("H26").Value = 0
$ws.Range("H27").Value = 0
$ws.Range("H28").Value = 0
$ws.Range("H29").Value = 0
$ws.Range("H30").Value = 0
$ws.Range("H31").Value = 0
$ws.Range("H32").Value = 1
